# trace length matching & generate gerber files
# Adds a "Via 0.45" length column (E / J) next to the existing "Length"
# columns, adds a BA0 pin (replacing BA2) to the sorted pin table, adds a
# new pad-to-pad measurement block (D26:D29) and restructures the
# MIN/MAX/AVG/DIFF summary block for the second (sorted) table to include
# AVG and the new "Via 0.45" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# New "Via 0.45" column for the first (unsorted) table: E5:E17
# ---------------------------------------------------------------------
$ws.Cells.Item(5, 5).Value = "Via 0.45"

$viaFirst = @(
    @(6, 75),
    @(7, 74.84),
    @(8, 75.018000000000001),
    @(9, 74.986000000000004),
    @(10, 74.998999999999995),
    @(11, 75)
)
foreach ($pair in $viaFirst) {
    $ws.Cells.Item($pair[0], 5).Value = $pair[1]
}

$ws.Cells.Item(14, 5).Formula = "=AVERAGE(E6:E11)"
$ws.Cells.Item(15, 5).Formula = "=MIN(E6:E11)"
$ws.Cells.Item(16, 5).Formula = "=MAX(E6:E11)"
$ws.Cells.Item(17, 5).Formula = "=E16-E15"

# ---------------------------------------------------------------------
# Pin table (G6:G44): BA2 -> BA0 (re-sorts alphabetically to the top)
# ---------------------------------------------------------------------
$pins = @(
    @(6, "BA0"),
    @(7, "BA1"),
    @(8, "BL0"),
    @(9, "BL1"),
    @(10, "CKE"),
    @(11, "CLK"),
    @(12, "D00"),
    @(13, "D01"),
    @(14, "D02"),
    @(15, "D03"),
    @(16, "D04"),
    @(17, "D05"),
    @(18, "D06"),
    @(19, "D07"),
    @(20, "D08"),
    @(21, "D09"),
    @(22, "D10"),
    @(23, "D11"),
    @(24, "D12"),
    @(25, "D13"),
    @(26, "D14"),
    @(27, "D15"),
    @(28, "FMC_A00"),
    @(29, "FMC_A01"),
    @(30, "FMC_A02"),
    @(31, "FMC_A03"),
    @(32, "FMC_A04"),
    @(33, "FMC_A05"),
    @(34, "FMC_A06"),
    @(35, "FMC_A07"),
    @(36, "FMC_A08"),
    @(37, "FMC_A09"),
    @(38, "FMC_A10"),
    @(39, "FMC_A11"),
    @(40, "FMC_A12"),
    @(41, "nCAS"),
    @(42, "nE1(nCS)"),
    @(43, "nRAS"),
    @(44, "nWE")
)
foreach ($pair in $pins) {
    $ws.Cells.Item($pair[0], 7).Value = $pair[1]
}

# ---------------------------------------------------------------------
# New "Via 0.45" column for the sorted table: J5:J44, J47:J51
# ---------------------------------------------------------------------
$ws.Cells.Item(5, 10).Value = "Via 0.45"

for ($r = 6; $r -le 44; $r++) {
    $ws.Cells.Item($r, 10).Value = 82.241
}
$ws.Cells.Item(17, 10).Value = 82.242000000000004

# ---------------------------------------------------------------------
# New pad-to-pad measurement block: D26:D29
# ---------------------------------------------------------------------
$ws.Cells.Item(26, 4).Value = 80.971000000000004
$ws.Cells.Item(27, 4).Value = 80.207999999999998
$ws.Cells.Item(28, 4).Formula = "=D26-D27"
$ws.Cells.Item(29, 4).Formula = "=82.241-D28"

# ---------------------------------------------------------------------
# Restructured summary block: MAX / MIN / AVG / DIFF (rows 47-51)
# ---------------------------------------------------------------------
$ws.Cells.Item(47, 8).Value = "MAX"
$ws.Cells.Item(47, 9).Formula = "=MAX(I6:I44)"
$ws.Cells.Item(47, 10).Formula = "=MAX(J6:J44)"

$ws.Cells.Item(48, 8).Value = "MIN"
$ws.Cells.Item(48, 9).Formula = "=MIN(I6:I44)"
$ws.Cells.Item(48, 10).Formula = "=MIN(J6:J44)"

$ws.Cells.Item(49, 8).Value = "AVG"
$ws.Cells.Item(49, 9).Formula = "=AVERAGE(I6:I44)"
$ws.Cells.Item(49, 10).Formula = "=AVERAGE(J6:J44)"
$ws.Cells.Item(49, 9).NumberFormat = "0.000"
$ws.Cells.Item(49, 10).NumberFormat = "0.000"

$ws.Cells.Item(51, 8).Value = "DIFF"
$ws.Cells.Item(51, 9).Formula = "=I47-I48"
$ws.Cells.Item(51, 10).Formula = "=J47-J48"

# ---------------------------------------------------------------------
# Column I width
# ---------------------------------------------------------------------
$ws.Columns.Item(9).ColumnWidth = 11.3

# ---------------------------------------------------------------------
# Selection
# ---------------------------------------------------------------------
$ws.Range("F25").Select()
